# Append one new telemetry row to each of the four worksheets, matching
# the "ID" log rows already present at the bottom of each sheet.
#
# Columns: A=time (date serial), B=总长, C=ID, D=实际长度, E=和校验,
#          F=总长_DEC, G=ID_DEC, H=实际长度_DEC, I=和校验_DEC

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param($Sheet, $Row, $TimeValue, $B, $C, $D, $E, $F, $G, $GIsText, $H, $I)

    # Column A: keep the same date/time display format as the row above it.
    $Sheet.Cells.Item($Row, 1).NumberFormat = $Sheet.Cells.Item($Row - 1, 1).NumberFormat
    $Sheet.Cells.Item($Row, 1).Value = $TimeValue

    $Sheet.Cells.Item($Row, 2).Value = $B
    $Sheet.Cells.Item($Row, 3).Value = $C
    $Sheet.Cells.Item($Row, 4).Value = $D
    $Sheet.Cells.Item($Row, 5).Value = $E

    $Sheet.Cells.Item($Row, 6).Value = $F

    if ($GIsText -eq 1) {
        # Huge integer that must stay text (would otherwise lose precision
        # as a double), exactly like the existing row in this sheet.
        $Sheet.Cells.Item($Row, 7).NumberFormat = "@"
        $Sheet.Cells.Item($Row, 7).Value = $G
    } else {
        $Sheet.Cells.Item($Row, 7).Value = $G
    }

    $Sheet.Cells.Item($Row, 8).Value = $H
    $Sheet.Cells.Item($Row, 9).Value = $I
}

# ROW50-FE-LIFTER: add row 45
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
$t1 = [double]"45746.68190453704"
$g1 = [double]"5.68631262647114E+23"
Add-LogRow $ws1 45 $t1 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x6a" "0xe" 400 $g1 0 362 14

# ROW50-MID-LIFTER: add row 47
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
$t2 = [double]"45746.6508912037"
$g2 = "568631262647113771663628"
Add-LogRow $ws2 47 $t2 "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," "0x01,0x6e" "0x19" 400 $g2 1 366 25

# ROW11-FE-LIFTER: add row 45
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
$t3 = [double]"45746.70389655093"
$g3 = [double]"5.68631262647114E+23"
Add-LogRow $ws3 45 $t3 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," "0x01,0x6a" "0x14" 400 $g3 0 362 20

# ROW11-MID-LIFTER: add row 45
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
$t4 = [double]"45746.84336512732"
$g4 = [double]"5.68631262647114E+23"
Add-LogRow $ws4 45 $t4 "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," "0x01,0x72" "0x19" 400 $g4 0 370 25
